$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1): rename + insert new Std_* columns, shift Obs_Prob to O ---
$ws.Range("A1").Value = "#_Agents"
$ws.Range("B1").Value = "Coverage"
$ws.Range("C1").Value = "Avg_Total_Rounds"
$ws.Range("D1").Value = "Avg_Expl_Cost"
$ws.Range("E1").Value = "Avg_Expl_Eff"
$ws.Range("F1").Value = "Avg_Round_Time"
$ws.Range("G1").Value = "Avg_Agent_Step_Time"
$ws.Range("H1").Value = "Avg_Experiment_Time"
$ws.Range("I1").Value = "Std_Total_Rounds"
$ws.Range("J1").Value = "Std_Expl_Cost"
$ws.Range("K1").Value = "Std_Expl_Eff"
$ws.Range("L1").Value = "Std_Round_Time"
$ws.Range("M1").Value = "Std_Agent_Step_Time"
$ws.Range("N1").Value = "Std_Experiment_Time"
$ws.Range("O1").Value = "Obs_Prob"

# New header cells K1:O1 need the same bold/centered/bordered style as the rest of row 1 (style copied from A1)
$ws.Range("A1").Copy()
$ws.Range("K1:O1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update data rows (rows 2-13) with the recomputed values (rows 2/3, 4/5, ... swapped by Obs_Prob, plus new Std_* stats) ---
# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 54.634
$ws.Range("D2").Value = 54.634
$ws.Range("E2").Value = 3.13345588
$ws.Range("F2").Value = 0.1216005
$ws.Range("G2").Value = 0.1216005
$ws.Range("H2").Value = 6.55444244
$ws.Range("I2").Value = 6.329943156199859
$ws.Range("J2").Value = 6.329943156199859
$ws.Range("K2").Value = 0.3522844533788149
$ws.Range("L2").Value = 0.01673302632191139
$ws.Range("M2").Value = 0.01673302632191139
$ws.Range("N2").Value = 0.4966816585414741
$ws.Range("O2").Value = 0.15

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 90.928
$ws.Range("D3").Value = 90.928
$ws.Range("E3").Value = 1.89545856
$ws.Range("F3").Value = 0.07737415999999998
$ws.Range("G3").Value = 0.07737415999999998
$ws.Range("H3").Value = 6.93212752
$ws.Range("I3").Value = 12.88871809406259
$ws.Range("J3").Value = 12.88871809406259
$ws.Range("K3").Value = 0.2652280383235383
$ws.Range("L3").Value = 0.01093635110287204
$ws.Range("M3").Value = 0.01093635110287204
$ws.Range("N3").Value = 0.6891824498739635
$ws.Range("O3").Value = 0.85

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 31.598
$ws.Range("D4").Value = 63.16
$ws.Range("E4").Value = 2.78779306
$ws.Range("F4").Value = 0.17614982
$ws.Range("G4").Value = 0.08807480000000001
$ws.Range("H4").Value = 2.69610944
$ws.Range("I4").Value = 6.452260016353748
$ws.Range("J4").Value = 12.90206675898503
$ws.Range("K4").Value = 0.5647982108374342
$ws.Range("L4").Value = 0.03939822215768082
$ws.Range("M4").Value = 0.0196996031149303
$ws.Range("N4").Value = 0.4394729316909797
$ws.Range("O4").Value = 0.15

# Row 5
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 48.862
$ws.Range("D5").Value = 95.89
$ws.Range("E5").Value = 1.81214186
$ws.Range("F5").Value = 0.11770542
$ws.Range("G5").Value = 0.05885248
$ws.Range("H5").Value = 2.80656324
$ws.Range("I5").Value = 8.693567289924712
$ws.Range("J5").Value = 16.08948859617897
$ws.Range("K5").Value = 0.3035120056377885
$ws.Range("L5").Value = 0.02206215434022398
$ws.Range("M5").Value = 0.01103121402272798
$ws.Range("N5").Value = 0.3700467820198805
$ws.Range("O5").Value = 0.85

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 16.002
$ws.Range("D6").Value = 63.93
$ws.Range("E6").Value = 2.8242366
$ws.Range("F6").Value = 0.2400204
$ws.Range("G6").Value = 0.06000524
$ws.Range("H6").Value = 0.93188626
$ws.Range("I6").Value = 4.182461232806331
$ws.Range("J6").Value = 16.69726616529753
$ws.Range("K6").Value = 0.7267734663607719
$ws.Range("L6").Value = 0.06102469562063528
$ws.Range("M6").Value = 0.01525627562790044
$ws.Range("N6").Value = 0.2470585752379577
$ws.Range("O6").Value = 0.15

# Row 7
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 25.7
$ws.Range("D7").Value = 94.71
$ws.Range("E7").Value = 1.84282608
$ws.Range("F7").Value = 0.17455174
$ws.Range("G7").Value = 0.04363776
$ws.Range("H7").Value = 1.09179158
$ws.Range("I7").Value = 5.766498505160652
$ws.Range("J7").Value = 16.76962553789925
$ws.Range("K7").Value = 0.3413453344326115
$ws.Range("L7").Value = 0.03912838032766668
$ws.Range("M7").Value = 0.009782105138685056
$ws.Range("N7").Value = 0.2378935879793884
$ws.Range("O7").Value = 0.85

# Row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 10.184
$ws.Range("D8").Value = 61.018
$ws.Range("E8").Value = 2.98890054
$ws.Range("F8").Value = 0.2835155999999999
$ws.Range("G8").Value = 0.0472524
$ws.Range("H8").Value = 0.4683868
$ws.Range("I8").Value = 2.999355976296427
$ws.Range("J8").Value = 17.99603819682441
$ws.Range("K8").Value = 0.8027512641932345
$ws.Range("L8").Value = 0.08466702085393439
$ws.Range("M8").Value = 0.01411132435990261
$ws.Range("N8").Value = 0.161982069523309
$ws.Range("O8").Value = 0.15

# Row 9
$ws.Range("A9").Value = 6
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 17.51
$ws.Range("D9").Value = 89.138
$ws.Range("E9").Value = 1.95902358
$ws.Range("F9").Value = 0.19246968
$ws.Range("G9").Value = 0.0320784
$ws.Range("H9").Value = 0.5475926799999999
$ws.Range("I9").Value = 4.557130113321264
$ws.Range("J9").Value = 16.18873392864386
$ws.Range("K9").Value = 0.362106208537089
$ws.Range("L9").Value = 0.04751380519033315
$ws.Range("M9").Value = 0.007919173521177386
$ws.Range("N9").Value = 0.1531297898094024
$ws.Range("O9").Value = 0.85

# Row 10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 7.046
$ws.Range("D10").Value = 56.166
$ws.Range("E10").Value = 3.20673472
$ws.Range("F10").Value = 0.28265714
$ws.Range("G10").Value = 0.03533212
$ws.Range("H10").Value = 0.24888046
$ws.Range("I10").Value = 1.796200621686423
$ws.Range("J10").Value = 14.27325806505239
$ws.Range("K10").Value = 0.8236979712140158
$ws.Range("L10").Value = 0.08970653709298353
$ws.Range("M10").Value = 0.01121325574998757
$ws.Range("N10").Value = 0.1057978885897469
$ws.Range("O10").Value = 0.15

# Row 11
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 13.398
$ws.Range("D11").Value = 83.142
$ws.Range("E11").Value = 2.1007238
$ws.Range("F11").Value = 0.19063272
$ws.Range("G11").Value = 0.02382906
$ws.Range("H11").Value = 0.31358678
$ws.Range("I11").Value = 3.657643053307802
$ws.Range("J11").Value = 15.35869835729776
$ws.Range("K11").Value = 0.3799389611212778
$ws.Range("L11").Value = 0.05346330580066923
$ws.Range("M11").Value = 0.006683218606201398
$ws.Range("N11").Value = 0.1076538559906281
$ws.Range("O11").Value = 0.85

# Row 12
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 5.858
$ws.Range("D12").Value = 58.314
$ws.Range("E12").Value = 3.10643778
$ws.Range("F12").Value = 0.2789528
$ws.Range("G12").Value = 0.02789532
$ws.Range("H12").Value = 0.16727792
$ws.Range("I12").Value = 1.584568621243921
$ws.Range("J12").Value = 15.73985037348065
$ws.Range("K12").Value = 0.8258367178945953
$ws.Range("L12").Value = 0.1026587928472652
$ws.Range("M12").Value = 0.01026589665669453
$ws.Range("N12").Value = 0.08858626298955448
$ws.Range("O12").Value = 0.15

# Row 13
$ws.Range("A13").Value = 10
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 10.966
$ws.Range("D13").Value = 76.74
$ws.Range("E13").Value = 2.29169176
$ws.Range("F13").Value = 0.16752428
$ws.Range("G13").Value = 0.01675238
$ws.Range("H13").Value = 0.17951754
$ws.Range("I13").Value = 3.438093217168476
$ws.Range("J13").Value = 15.54747078130044
$ws.Range("K13").Value = 0.4587957001440224
$ws.Range("L13").Value = 0.0479873071068625
$ws.Range("M13").Value = 0.004798775283762316
$ws.Range("N13").Value = 0.06813145271786956
$ws.Range("O13").Value = 0.85
